$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.484.65"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "1.650.88"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").Value = "'300.37"
$ws.Range("E6").Value = "  -0.70%  "

$ws.Range("D7").Value = "'0.3785"
$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("D8").Value = "'50.67"
$ws.Range("E8").Value = "  -0.90%  "

$ws.Range("D9").Value = "'0.3500"
$ws.Range("E9").Value = "  -2.58%  "

$ws.Range("D10").Value = "'1.225"
$ws.Range("E10").Value = "  -1.31%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "'22.09"
$ws.Range("E13").Value = "  -1.60%  "

$ws.Range("D14").Value = "'6.316"
$ws.Range("E14").Value = "  -2.66%  "

$ws.Range("D15").Value = "'7.272"
$ws.Range("E15").Value = "  -2.99%  "

$ws.Range("D16").Value = "'0.00001211"
$ws.Range("E16").Value = "  -0.86%  "

$ws.Range("D17").Value = "1.648.69"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").Value = "'95.39"
$ws.Range("E18").Value = "  -2.14%  "

$ws.Range("D19").Value = "'0.06969"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").Value = "'6.638"
$ws.Range("E20").Value = "  -2.77%  "

$ws.Range("E21").Value = "  -1.20%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("E23").Value = "  -1.59%  "

$ws.Range("D24").Value = "23.497.64"
$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("D25").Value = "'2.424"
$ws.Range("E25").Value = "  -3.11%  "

$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("D27").Value = "'21.10"
$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").Value = "'151.83"
$ws.Range("E28").Value = "  -0.67%  "

$ws.Range("D29").Value = "'5.188"
$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("D30").Value = "'131.87"
$ws.Range("E30").Value = "  -1.36%  "

$ws.Range("D31").Value = "1.834.35"
$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("D32").Value = "'6.924"
$ws.Range("E32").Value = "  -3.94%  "

$ws.Range("D33").Value = "'2.143"
$ws.Range("E33").Value = "  -4.81%  "

$ws.Range("D34").Value = "'11.19"
$ws.Range("E34").Value = "  -7.43%  "

$ws.Range("D35").Value = "'0.9892"
$ws.Range("E35").Value = "  -6.58%  "

$ws.Range("D36").Value = "'0.02728"
$ws.Range("E36").Value = "  -2.71%  "

$ws.Range("D37").Value = "'0.08767"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("E38").Value = "  -2.99%  "

$ws.Range("D39").Value = "'0.2425"
$ws.Range("E39").Value = "  -3.04%  "

$ws.Range("D40").Value = "'0.06837"
$ws.Range("E40").Value = "  -2.71%  "

$ws.Range("D41").Value = "'12.87"
$ws.Range("E41").Value = "  -2.57%  "

$ws.Range("D42").Value = "'0.6918"
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("D43").Value = "'1.298"
$ws.Range("E43").Value = "  -2.90%  "

$ws.Range("D44").Value = "'15.65"
$ws.Range("E44").Value = "  -2.25%  "

$ws.Range("D45").Value = "'0.9991"
$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("D46").Value = "'0.6395"
$ws.Range("E46").Value = "  -2.16%  "

$ws.Range("D47").Value = "'2.253"
$ws.Range("E47").Value = "  -2.48%  "

$ws.Range("D48").Value = "'3.924"
$ws.Range("E48").Value = "  -0.88%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.07683"
$ws.Range("E49").Value = "  -2.93%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'127.14"
$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("D51").Value = "'1.242"
$ws.Range("E51").Value = "  +3.23%  "
